$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 91820748
$ws.Range("B2").Value = 90647
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 4362
$ws.Range("F2").Value = 'Blå taggsvamp'
$ws.Range("G2").Value = 'Hydnellum caeruleum'
$ws.Range("H2").Value = '(Hornem.) P.Karst.'
$ws.Range("Q2").Value = 655526.8862207049
$ws.Range("R2").Value = 7353400.057194735

# Row 3
$ws.Range("A3").Value = 91820740
$ws.Range("B3").Value = 89633
$ws.Range("D3").Value = 'VU'
$ws.Range("E3").Value = 65
$ws.Range("F3").Value = 'Fläckporing'
$ws.Range("G3").Value = 'Anthoporia albobrunnea'
$ws.Range("H3").Value = '(Romell) Karasiński & Niemelä'
$ws.Range("Q3").Value = 655588.8241655316
$ws.Range("R3").Value = 7352968.129340165

# Row 4
$ws.Range("A4").Value = 91820735
$ws.Range("B4").Value = 90665
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 4366
$ws.Range("F4").Value = 'Skarp dropptaggsvamp'
$ws.Range("G4").Value = 'Hydnellum peckii'
$ws.Range("H4").Value = 'Banker'
$ws.Range("Q4").Value = 655685.1562783264
$ws.Range("R4").Value = 7352929.848888704

# Row 5
$ws.Range("A5").Value = 91820754
$ws.Range("B5").Value = 90665
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 4366
$ws.Range("F5").Value = 'Skarp dropptaggsvamp'
$ws.Range("G5").Value = 'Hydnellum peckii'
$ws.Range("H5").Value = 'Banker'
$ws.Range("Q5").Value = 655600.9169711781
$ws.Range("R5").Value = 7352968.799517729
